$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.142.17'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.67%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.470.87'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.63%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '560.55'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.64%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '162.80'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.18%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.507'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.14%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.469.07'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.58%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.151'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.61%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.61%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.88'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.16%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.332'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -3.41%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '69.002.52'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.76%  '
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.06%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000169'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -1.53%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '23.66'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.42%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.477.07'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.05%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.75'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.49%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '338.14'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -3.26%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.99'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -2.92%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.80'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.61%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.35%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.05%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '67.19'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.79%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.68'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -2.73%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.601.03'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.70%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.30'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.09%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.73%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0821'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -2.83%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.20'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -2.12%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.04%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '431.86'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.11%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.14'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -3.44%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -3.50%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '156.83'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.00%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.22%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.01%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.109'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.78%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '17.81'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.91%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.301'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -1.97%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.43'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -2.29%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.47'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -4.76%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.08'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -1.98%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.07'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.47%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '132.70'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -2.13%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.62%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.32%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.485'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.23%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.561'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.63%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0916'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.15%  '
